# "Generate Report for Handoff"
# Updates the localization-status report: flips the in-progress rows from
# "In Translation" to "Ready for handoff" and refreshes the report's
# generation timestamps, then widens the Status columns so the longer
# "Ready for handoff" text fits (as the report generator's column
# auto-sizing would do).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-11-08 23:31:36"

# Widen the status columns (zh-cn / de-de) to fit the new text.
$wsOverview.Columns("E:F").ColumnWidth = 16.333333333333332

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-11-08 23:31:23"
$wsZhCn.Columns("C").ColumnWidth = 16.333333333333332

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-11-08 23:31:36"
$wsDeDe.Columns("C").ColumnWidth = 16.333333333333332
